$d = $word.ActiveDocument

# Locate the paragraph that contains "Ver no Jupiter Salvar em pdf Salvar em docx".
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $prev = $target.Previous()
    $next = $target.Next()

    # Delete the three paragraphs individually (each paragraph's own range,
    # including its paragraph mark) in reverse order so earlier deletions
    # don't shift the ranges of the ones still pending.
    $next.Range.Delete()
    $target.Range.Delete()
    $prev.Range.Delete()
}
